# Fix formatting of the country-paragraph text in the "overview" sheet.
# For each country paragraph, insert a trailing "<br><br>" right after the
# final sentence of the paragraph (i.e. right before the next country's
# bolded heading, or before the closing "Faroe Islands" sentence for the
# last country paragraph).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cell = $ws.Range("B4")

$text = $cell.Value2

$pairs = @(
    @{ Old = "common cuttlefish, squid, rays, and cod.`n<b>Denmark</b><br>"; New = "common cuttlefish, squid, rays, and cod.<br><br>`n<b>Denmark</b><br>" },
    @{ Old = "for the period.`n<b>France</b><br>"; New = "for the period.<br><br>`n<b>France</b><br>" },
    @{ Old = "and horse-mackerel.`n<B>Germany</b><br>"; New = "and horse-mackerel.<br><br>`n<B>Germany</b><br>" },
    @{ Old = "for the last ten years.`n<b>Netherlands</b><br>"; New = "for the last ten years.<br><br>`n<b>Netherlands</b><br>" },
    @{ Old = "mackerel, and horse mackerel.`n<b>Norway</b><br>"; New = "mackerel, and horse mackerel.<br><br>`n<b>Norway</b><br>" },
    @{ Old = "of the total landings (up from 75% in ten years; +40% in absolute landing tonnage).`n<b>Sweden</b><br>"; New = "of the total landings (up from 75% in ten years; +40% in absolute landing tonnage).<br><br>`n<b>Sweden</b><br>" },
    @{ Old = "number of vessels over the last decade.`n<b>UK (England)</b><br>"; New = "number of vessels over the last decade.<br><br>`n<b>UK (England)</b><br>" },
    @{ Old = "unchanged since 2009. This largest category of vessel targets finfish, mostly demersal but some pelagic.`n<b>UK (Scotland)</b><br>"; New = "unchanged since 2009. This largest category of vessel targets finfish, mostly demersal but some pelagic.<br><br>`n<b>UK (Scotland)</b><br>" },
    @{ Old = "pelagic species are harvested by 18 large vessels, primarily using pelagic trawls.`nThe Faroe Islands"; New = "pelagic species are harvested by 18 large vessels, primarily using pelagic trawls.<br><br>`nThe Faroe Islands" }
)

foreach ($pair in $pairs) {
    $old = $pair.Old
    $new = $pair.New
    if ($text.IndexOf($old) -lt 0) {
        throw "Anchor not found: $old"
    }
    $text = $text.Replace($old, $new)
}

$cell.Value2 = $text
